$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1117
$ws1.Range("F4").Value = 1799
$ws1.Range("F6").Value = 364
$ws1.Range("F7").Value = 214

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1117
$ws4.Range("F4").Value = 1799
$ws4.Range("F7").Value = 364
$ws4.Range("F8").Value = 214
